# Adds the new 2017-05-11 reimbursement data to each relevant sheet of the
# "莆田项目部公共开支明细表" workbook, mirroring the rows that a user would
# have typed directly into Excel for that date.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# Sheet 1 - 逐日消费统计表 (daily consumption statistics)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Cells.Item(5, 1).Value = 4
$ws1.Cells.Item(5, 2).Value = 20170511
$ws1.Cells.Item(5, 3).Formula = "=3+35+100+4+51+5"
$ws1.Cells.Item(5, 4).Value = "林迪南、吴紫东、郑景祥"

# Column E on row 5 needs the same style as the rows above it (s=1).
$ws1.Range("E4").Copy() | Out-Null
$ws1.Range("E5").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false
$ws1.Cells.Item(5, 5).Value = "林迪南"

# ---------------------------------------------------------------------
# Sheet 2 - 资金垫付表 (fund advance table)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Cells.Item(4, 1).Value = 3
$ws2.Cells.Item(4, 2).Value = 20170511

# Column C keeps the bold-ish "垫付者" style used on rows 2-3 (s=1).
$ws2.Range("C2").Copy() | Out-Null
$ws2.Range("C4").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false
$ws2.Cells.Item(4, 3).Value = "林迪南"

$ws2.Cells.Item(4, 4).Value = 198

# ---------------------------------------------------------------------
# Sheet 3 - 人员固定消费表 (personnel fixed consumption table)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Cells.Item(4, 1).Value = 3
$ws3.Cells.Item(4, 2).Value = 20170511
$ws3.Cells.Item(4, 3).Value = "郑景祥"
$ws3.Cells.Item(4, 4).Value = "伙食"
$ws3.Cells.Item(4, 5).Formula = "=100/3+51/2"

$ws3.Cells.Item(5, 1).Value = 4
$ws3.Cells.Item(5, 2).Value = 20170511
$ws3.Cells.Item(5, 3).Value = "郑景祥"
$ws3.Cells.Item(5, 4).Value = "公交车"
$ws3.Cells.Item(5, 5).Value = 5

$ws3.Cells.Item(6, 1).Value = 5
$ws3.Cells.Item(6, 2).Value = 20170511
$ws3.Cells.Item(6, 3).Value = "林迪南"
$ws3.Cells.Item(6, 4).Value = "伙食"
$ws3.Cells.Item(6, 5).Formula = "=100/3+51/2+4"

$ws3.Cells.Item(7, 1).Value = 6
$ws3.Cells.Item(7, 2).Value = 20170511
$ws3.Cells.Item(7, 3).Value = "吴紫东"
$ws3.Cells.Item(7, 4).Value = "伙食"
$ws3.Cells.Item(7, 5).Formula = "=100/3"

# ---------------------------------------------------------------------
# Sheet 4 - 公共消费权重表 (public consumption weight table)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Cells.Item(8, 1).Value = 7
$ws4.Cells.Item(8, 2).Value = 20170511
$ws4.Cells.Item(8, 3).Value = "林迪南"
$ws4.Cells.Item(8, 4).Value = 0

$ws4.Cells.Item(9, 2).Value = 20170511
$ws4.Cells.Item(9, 3).Value = "郑景祥"
$ws4.Cells.Item(9, 4).Value = 0

$ws4.Cells.Item(10, 2).Value = 20170511
$ws4.Cells.Item(10, 3).Value = "吴紫东"
$ws4.Cells.Item(10, 4).Value = 0

# ---------------------------------------------------------------------
# Sheet 5 - 莆田维养项目部签到表 (sign-in sheet)
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

$ws5.Cells.Item(6, 2).Value = "郑景祥"
$ws5.Cells.Item(6, 3).Value = 20170511

# D6 carries the same "s=6" style used by the D column elsewhere.
$ws5.Range("D2").Copy() | Out-Null
$ws5.Range("D6").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false
$ws5.Cells.Item(6, 4).Value = "莆田"

$ws5.Cells.Item(6, 5).Value = "莆田"

$ws5.Cells.Item(7, 2).Value = "王燊"
$ws5.Cells.Item(7, 3).Value = 20170511
$ws5.Cells.Item(7, 4).Value = "福州"
$ws5.Cells.Item(7, 5).Value = "福州"

# ---------------------------------------------------------------------
# Sheet 6 - 莆田绶溪公园施工监控项目签到表 (sign-in sheet)
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

$ws6.Cells.Item(6, 2).Value = "黄学漾"
$ws6.Cells.Item(6, 3).Value = 20170511
$ws6.Cells.Item(6, 4).Value = "福州"
$ws6.Cells.Item(6, 5).Value = "福州"

$ws6.Cells.Item(7, 2).Value = "王兆林"
$ws6.Cells.Item(7, 3).Value = 20170511

# ---------------------------------------------------------------------
# Restore the selections/active cells shown in the final workbook, ending
# on sheet 3 (人员固定消费表) so it becomes the active tab, as in the diff.
# ---------------------------------------------------------------------
$ws1.Range("C5").Select() | Out-Null
$ws2.Range("B5").Select() | Out-Null
$ws4.Range("D11").Select() | Out-Null
$ws5.Range("C8").Select() | Out-Null
$ws6.Range("B1").Select() | Out-Null
$ws3.Range("E7").Select() | Out-Null
